# Generate Report for Handback
#
# This reproduces, via the Excel object model, the "handback" update that
# fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns for the zh-cn and de-de sheets, and
# flips the overall Status from "In Translation" to
# "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/e70270d8d66483fe91af77ee81849217d100b046/e2e"

$newStatus = "Handed back: in sync with en-US"

# Source file rows (same 3 content rows on every per-locale sheet).
$rows = @(
    @{ Row = 2; Name = "3b4f7f11-227e-4272-950b-c531c7b6c03a.md" },
    @{ Row = 3; Name = "cf5fc51b-3d98-4493-9d84-8cbb994eacd0.md" },
    @{ Row = 4; Name = "toc.md" }
)

# ---------------------------------------------------------------------
# 1. Overview sheet: Status column text (zh-cn => E, de-de => F).
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $ov.Range("E" + $r.Row).Value = $newStatus
    $ov.Range("F" + $r.Row).Value = $newStatus
}

# ---------------------------------------------------------------------
# 2. zh-cn sheet: Status + Latest Target File / Handback File / DateTime.
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zhHandbackFiles = @(
    "3b4f7f11-227e-4272-950b-c531c7b6c03a.66cf05ae6d1b6e70cb496d0ee8bba94f91211062.zh-cn.xlf",
    "cf5fc51b-3d98-4493-9d84-8cbb994eacd0.ec7a02cb049541d09450274516252969c3ccc104.zh-cn.xlf",
    "toc.052a173bd0b736745b4800c3b8aeca39fe30b2dd.zh-cn.xlf"
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $row = $rows[$i].Row
    $name = $rows[$i].Name

    # Status
    $zh.Range("C" + $row).Value = $newStatus

    # Latest Target File (J) - same source file, now hyperlinked.
    $jCell = $zh.Range("J" + $row)
    $jCell.Value = $name
    $zh.Hyperlinks.Add($jCell, $baseUrl + "/" + $name, "", "", $name)

    # Latest Handback File (K)
    $zh.Range("K" + $row).Value = $zhHandbackFiles[$i]
}

# Latest Handback DateTime (L) - same instant for every row on this sheet.
$zh.Range("L2").Value = "2017-04-28 03:05:27"
$zh.Range("L3").Value = "2017-04-28 03:05:27"
$zh.Range("L4").Value = "2017-04-28 03:05:27"

# ---------------------------------------------------------------------
# 3. de-de sheet: Status + Latest Target File / Handback File / DateTime.
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$deHandbackFiles = @(
    "3b4f7f11-227e-4272-950b-c531c7b6c03a.66cf05ae6d1b6e70cb496d0ee8bba94f91211062.de-de.xlf",
    "cf5fc51b-3d98-4493-9d84-8cbb994eacd0.ec7a02cb049541d09450274516252969c3ccc104.de-de.xlf",
    "toc.052a173bd0b736745b4800c3b8aeca39fe30b2dd.de-de.xlf"
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $row = $rows[$i].Row
    $name = $rows[$i].Name

    # Status
    $de.Range("C" + $row).Value = $newStatus

    # Latest Target File (J) - same source file, now hyperlinked.
    $jCell = $de.Range("J" + $row)
    $jCell.Value = $name
    $de.Hyperlinks.Add($jCell, $baseUrl + "/" + $name, "", "", $name)

    # Latest Handback File (K)
    $de.Range("K" + $row).Value = $deHandbackFiles[$i]
}

# Latest Handback DateTime (L) - de-de was handed back a bit later than zh-cn.
$de.Range("L2").Value = "2017-04-28 03:05:59"
$de.Range("L3").Value = "2017-04-28 03:05:59"
$de.Range("L4").Value = "2017-04-28 03:05:59"

# ---------------------------------------------------------------------
# 4. Column width touch-up: widen columns that now hold the long
#    "Handed back: in sync with en-US" / file-name strings, mirroring the
#    autofit Excel performs when a column's content grows.
# ---------------------------------------------------------------------
$ov.Columns.Item(5).ColumnWidth = 29.09
$ov.Columns.Item(6).ColumnWidth = 29.09

foreach ($ws in @($zh, $de)) {
    $ws.Columns.Item(3).ColumnWidth = 29.09
    $ws.Columns.Item(10).ColumnWidth = 39.17
    $ws.Columns.Item(11).ColumnWidth = 39.17
}
